$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Effort (D5) and Remain (E5) for the "Save items locally in json format" task
$ws.Range("D5").Value = 3.5
$ws.Range("E5").Value = 0.5

# Update the active selection to E6
$ws.Range("E6").Select()
